$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44365
$ws.Range("J2").Value = 580
$ws.Range("K2").Value = 20000
$ws.Range("M2").Value = 21103
$ws.Range("P2").Value = 1407

# Row 3
$ws.Range("D3").Value = 44406
$ws.Range("J3").Value = 400
$ws.Range("L3").Value = 22000
$ws.Range("M3").Value = 20850
$ws.Range("P3").Value = 1390

# Row 5
$ws.Range("D5").Value = 44483
$ws.Range("J5").Value = 220
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 18909
$ws.Range("P5").Value = 1261

# Row 6
$ws.Range("D6").Value = 44398
$ws.Range("J6").Value = 130
$ws.Range("L6").Value = 20000
$ws.Range("M6").Value = 20000
$ws.Range("P6").Value = 1333

# Row 7
$ws.Range("D7").Value = 44435
$ws.Range("J7").Value = 140
$ws.Range("K7").Value = 21000
$ws.Range("L7").Value = 23000
$ws.Range("M7").Value = 21714
$ws.Range("P7").Value = 1448

# Row 8
$ws.Range("D8").Value = 44476
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 22000
$ws.Range("M8").Value = 20909
$ws.Range("P8").Value = 1394

# Row 9
$ws.Range("D9").Value = 44400
$ws.Range("J9").Value = 130
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 24000
$ws.Range("M9").Value = 24000
$ws.Range("P9").Value = 1600

# Row 11
$ws.Range("D11").Value = 44396
$ws.Range("K11").Value = 22000
$ws.Range("L11").Value = 22000
$ws.Range("M11").Value = 22000
$ws.Range("P11").Value = 1467

# Row 12
$ws.Range("D12").Value = 44399
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 22000
$ws.Range("M12").Value = 22000
$ws.Range("P12").Value = 1467

# Row 13
$ws.Range("D13").Value = 44453
$ws.Range("J13").Value = 280
$ws.Range("K13").Value = 20000
$ws.Range("L13").Value = 22000
$ws.Range("M13").Value = 21286
$ws.Range("P13").Value = 1419

# Row 14
$ws.Range("D14").Value = 44449
$ws.Range("K14").Value = 22000
$ws.Range("L14").Value = 24000
$ws.Range("M14").Value = 23091
$ws.Range("P14").Value = 1539
